$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item('!!Compartment').Unprotect()
$wb.Worksheets.Item('!!Compound').Unprotect()
$wb.Worksheets.Item('!!Definition').Unprotect()
$wb.Worksheets.Item('!!Enzyme').Unprotect()
$wb.Worksheets.Item('!!FbcObjective').Unprotect()
$wb.Worksheets.Item('!!Gene').Unprotect()
$wb.Worksheets.Item('!!Layout').Unprotect()
$wb.Worksheets.Item('!!Measurement').Unprotect()
$wb.Worksheets.Item('!!PbConfig').Unprotect()
$wb.Worksheets.Item('!!Position').Unprotect()
$wb.Worksheets.Item('!!Protein').Unprotect()
$wb.Worksheets.Item('!!Quantity').Unprotect()
$wb.Worksheets.Item('!!QuantityInfo').Unprotect()
$wb.Worksheets.Item('!!QuantityMatrix').Unprotect()
$wb.Worksheets.Item('!!Reaction').Unprotect()
$wb.Worksheets.Item('!!ReactionStoichiometry').Unprotect()
$wb.Worksheets.Item('!!Regulator').Unprotect()
$wb.Worksheets.Item('!!Relation').Unprotect()
$wb.Worksheets.Item('!!Relationship').Unprotect()
$wb.Worksheets.Item('!!SparseMatrix').Unprotect()
$wb.Worksheets.Item('!!SparseMatrixColumn').Unprotect()
$wb.Worksheets.Item('!!SparseMatrixOrdered').Unprotect()
$wb.Worksheets.Item('!!SparseMatrixRow').Unprotect()
$wb.Worksheets.Item('!!StoichiometricMatrix').Unprotect()
$wb.Worksheets.Item('!!rxnconContingencyList').Unprotect()
$wb.Worksheets.Item('!!rxnconReactionList').Unprotect()

$wb.Worksheets.Item('!!Compartment').Range('A1').Value = '!!!ObjTables objTablesVersion=''0.0.8'' date=''2020-03-09 15:31:21'''
$wb.Worksheets.Item('!!Compartment').Range('A2').Value = '!!ObjTables type=''Data'' id=''Compartment'' name=''Compartment'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Compound').Range('A1').Value = '!!ObjTables type=''Data'' id=''Compound'' name=''Compound'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Definition').Range('A1').Value = '!!ObjTables type=''Data'' id=''Definition'' name=''Definition'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Enzyme').Range('A1').Value = '!!ObjTables type=''Data'' id=''Enzyme'' name=''Enzyme'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!FbcObjective').Range('A1').Value = '!!ObjTables type=''Data'' id=''FbcObjective'' name=''FbcObjective'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Gene').Range('A1').Value = '!!ObjTables type=''Data'' id=''Gene'' name=''Gene'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Layout').Range('A1').Value = '!!ObjTables type=''Data'' id=''Layout'' name=''Layout'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Measurement').Range('A1').Value = '!!ObjTables type=''Data'' id=''Measurement'' name=''Measurement'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!PbConfig').Range('A1').Value = '!!ObjTables type=''Data'' id=''PbConfig'' name=''PbConfig'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Position').Range('A1').Value = '!!ObjTables type=''Data'' id=''Position'' name=''Position'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Protein').Range('A1').Value = '!!ObjTables type=''Data'' id=''Protein'' name=''Protein'' date=''2020-03-09 15:31:21'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Quantity').Range('A1').Value = '!!ObjTables type=''Data'' id=''Quantity'' name=''Quantity'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!QuantityInfo').Range('A1').Value = '!!ObjTables type=''Data'' id=''QuantityInfo'' name=''QuantityInfo'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!QuantityMatrix').Range('A1').Value = '!!ObjTables type=''Data'' id=''QuantityMatrix'' name=''QuantityMatrix'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Reaction').Range('A1').Value = '!!ObjTables type=''Data'' id=''Reaction'' name=''Reaction'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!ReactionStoichiometry').Range('A1').Value = '!!ObjTables type=''Data'' id=''ReactionStoichiometry'' name=''ReactionStoichiometry'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Regulator').Range('A1').Value = '!!ObjTables type=''Data'' id=''Regulator'' name=''Regulator'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Relation').Range('A1').Value = '!!ObjTables type=''Data'' id=''Relation'' name=''Relation'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Relationship').Range('A1').Value = '!!ObjTables type=''Data'' id=''Relationship'' name=''Relationship'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!SparseMatrix').Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrix'' name=''SparseMatrix'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!SparseMatrixColumn').Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixColumn'' name=''SparseMatrixColumn'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!SparseMatrixOrdered').Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixOrdered'' name=''SparseMatrixOrdered'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!SparseMatrixRow').Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixRow'' name=''SparseMatrixRow'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!StoichiometricMatrix').Range('A1').Value = '!!ObjTables type=''Data'' id=''StoichiometricMatrix'' name=''StoichiometricMatrix'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!rxnconContingencyList').Range('A1').Value = '!!ObjTables type=''Data'' id=''rxnconContingencyList'' name=''rxnconContingencyList'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!rxnconReactionList').Range('A1').Value = '!!ObjTables type=''Data'' id=''rxnconReactionList'' name=''rxnconReactionList'' date=''2020-03-09 15:31:22'' objTablesVersion=''0.0.8'' tableFormat=''row'''

$wb.Worksheets.Item('!!Compartment').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Compound').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Definition').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Enzyme').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!FbcObjective').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Gene').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Layout').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Measurement').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!PbConfig').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Position').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Protein').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Quantity').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!QuantityInfo').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!QuantityMatrix').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Reaction').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!ReactionStoichiometry').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Regulator').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Relation').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!Relationship').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!SparseMatrix').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!SparseMatrixColumn').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!SparseMatrixOrdered').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!SparseMatrixRow').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!StoichiometricMatrix').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!rxnconContingencyList').Protect($null, $true, $true, $true)
$wb.Worksheets.Item('!!rxnconReactionList').Protect($null, $true, $true, $true)
